$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F ("想去人数") for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 33
$ws1.Range("F6").Value = 438
$ws1.Range("F8").Value = 1989
$ws1.Range("F10").Value = 35
$ws1.Range("F11").Value = 31
$ws1.Range("F12").Value = 1593
$ws1.Range("F13").Value = 1593
$ws1.Range("F14").Value = 1321
$ws1.Range("F15").Value = 52
$ws1.Range("F18").Value = 14
$ws1.Range("F20").Value = 442
$ws1.Range("F23").Value = 140
$ws1.Range("F24").Value = 7006
$ws1.Range("F25").Value = 7581
$ws1.Range("F26").Value = 35
$ws1.Range("F27").Value = 2
$ws1.Range("F28").Value = 177
$ws1.Range("F30").Value = 74
$ws1.Range("F31").Value = 216
$ws1.Range("F32").Value = 245
$ws1.Range("F37").Value = 1378
$ws1.Range("F38").Value = 14
$ws1.Range("F40").Value = 277
$ws1.Range("F45").Value = 212
$ws1.Range("F46").Value = 184
$ws1.Range("F47").Value = 79
$ws1.Range("F48").Value = 125

# Sheet "演出" (sheet2): update column F for several rows
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 23
$ws2.Range("F16").Value = 3
$ws2.Range("F17").Value = 281

# Sheet "本地生活" (sheet3): update column F for several rows
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2560
$ws3.Range("F5").Value = 117
$ws3.Range("F6").Value = 3

# Sheet "全部类型" (sheet4): update column F for several rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 117
$ws4.Range("F6").Value = 23
$ws4.Range("F7").Value = 438
$ws4.Range("F8").Value = 1990
$ws4.Range("F9").Value = 35
$ws4.Range("F10").Value = 31
$ws4.Range("F11").Value = 1593
$ws4.Range("F12").Value = 1593
$ws4.Range("F13").Value = 3
$ws4.Range("F14").Value = 1321
$ws4.Range("F15").Value = 52
$ws4.Range("F16").Value = 14
$ws4.Range("F18").Value = 442
$ws4.Range("F20").Value = 140
$ws4.Range("F22").Value = 7006
$ws4.Range("F23").Value = 7581
$ws4.Range("F24").Value = 35
$ws4.Range("F25").Value = 177
$ws4.Range("F26").Value = 74
$ws4.Range("F27").Value = 245
$ws4.Range("F30").Value = 14
$ws4.Range("F33").Value = 277
$ws4.Range("F43").Value = 212
$ws4.Range("F44").Value = 184
$ws4.Range("F45").Value = 79
$ws4.Range("F46").Value = 125
$ws4.Range("F48").Value = 3
$ws4.Range("F49").Value = 281

$wb.Save()
